$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2023-10-18 Wednesday", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "2023-10-19 Thursday", 2)

# Update each division problem in the table, cell by cell, addressing cells
# directly by (row, col) and assigning the new text to the cell's Range.
# (Some problem texts repeat, e.g. "89÷5=" appears twice but must map to
# different replacements, so each cell is targeted individually rather than
# relying on a document-wide Find/Replace.)
$t = $d.Tables.Item(1)

$cellEdits = @(
    @{ Row = 1;  Col = 1; New = "27÷5=" },
    @{ Row = 1;  Col = 2; New = "32÷4=" },
    @{ Row = 1;  Col = 3; New = "21÷4=" },
    @{ Row = 1;  Col = 4; New = "23÷3=" },
    @{ Row = 1;  Col = 5; New = "56÷8=" },

    @{ Row = 5;  Col = 1; New = "89÷2=" },
    @{ Row = 5;  Col = 2; New = "23÷2=" },
    @{ Row = 5;  Col = 3; New = "87÷8=" },
    @{ Row = 5;  Col = 4; New = "48÷6=" },
    @{ Row = 5;  Col = 5; New = "52÷2=" },

    @{ Row = 9;  Col = 1; New = "85÷7=" },
    @{ Row = 9;  Col = 2; New = "76÷8=" },
    @{ Row = 9;  Col = 3; New = "69÷4=" },
    @{ Row = 9;  Col = 4; New = "97÷5=" },
    @{ Row = 9;  Col = 5; New = "48÷3=" },

    @{ Row = 13; Col = 1; New = "70÷4=" },
    @{ Row = 13; Col = 2; New = "70÷4=" },
    @{ Row = 13; Col = 3; New = "29÷9=" },
    @{ Row = 13; Col = 4; New = "45÷6=" },
    @{ Row = 13; Col = 5; New = "21÷6=" },

    @{ Row = 17; Col = 1; New = "90÷4=" },
    @{ Row = 17; Col = 2; New = "46÷6=" },
    @{ Row = 17; Col = 3; New = "43÷4=" },
    @{ Row = 17; Col = 4; New = "38÷3=" },
    @{ Row = 17; Col = 5; New = "41÷7=" }
)

foreach ($edit in $cellEdits) {
    $cell = $t.Cell($edit.Row, $edit.Col)
    $cell.Range.Text = $edit.New
}
